$wb = $excel.ActiveWorkbook

$wsFP = $wb.Worksheets.Item("Faculties & Programs")
$wsFA = $wb.Worksheets.Item("Faculty & Academic")

# --- Fix typo in the "Faculty" (category) column of "Faculties & Programs":
#     "Varcity Clubs" -> "Varsity Clubs" (applies to every row using that
#     category, e.g. B208:B242) ---
$used = $wsFP.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $wsFP.Cells.Item($r, 2)
    if ($cell.Value2 -eq "Varcity Clubs") {
        $cell.Value = "Varsity Clubs"
    }
}

# --- Update view/selection state ---
# "Faculty & Academic" is no longer the active/selected tab; scroll it down
# and move its selection near the bottom of its data.
$wsFA.Activate()
$winFA = $excel.ActiveWindow
$winFA.ScrollRow = 172
$winFA.ScrollColumn = 1
$wsFA.Range("C241").Select()

# "Faculties & Programs" becomes the active/selected sheet, scrolled and
# selected near the bottom of its data.
$wsFP.Activate()
$winFP = $excel.ActiveWindow
$winFP.ScrollRow = 224
$winFP.ScrollColumn = 1
$wsFP.Range("E240").Select()
